$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyDesign")
# copy style from B6 (s=27, left/top only) onto B1 (originally s=27) - no-op test
$ws.Range("B6").Copy()
$ws.Range("B1:E2").PasteSpecial(-4122)  # xlPasteFormats = -4122
